# Easy-upload template: add 8 new "$vN" variable columns to Table1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$t = $ws.ListObjects.Item(1)

for ($i = 1; $i -le 8; $i++) {
    $null = $t.ListColumns.Add()
    $ws.Cells.Item(1, 3 + $i).Value = '$v' + $i
}
